$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 7: DAMSLTag "%" -> "sd", DialogAct "Uninterpretable" -> "Statement-non-opinion"
$ws.Range("I7").Value = "sd"
$ws.Range("J7").Value = "Statement-non-opinion"

# Row 26: DAMSLTag "sd" -> "sv", DialogAct "Statement-non-opinion" -> "Statement-opinion"
$ws.Range("I26").Value = "sv"
$ws.Range("J26").Value = "Statement-opinion"

# Row 30: DAMSLTag "sd" -> "sv", DialogAct "Statement-non-opinion" -> "Statement-opinion"
$ws.Range("I30").Value = "sv"
$ws.Range("J30").Value = "Statement-opinion"

# Row 31: DAMSLTag "aa" -> "sd", DialogAct "Agree/Accept" -> "Statement-non-opinion"
$ws.Range("I31").Value = "sd"
$ws.Range("J31").Value = "Statement-non-opinion"

# Row 35: DAMSLTag "sd" -> "sv", DialogAct "Statement-non-opinion" -> "Statement-opinion"
$ws.Range("I35").Value = "sv"
$ws.Range("J35").Value = "Statement-opinion"

# Row 36: DAMSLTag "sd" -> "sv", DialogAct "Statement-non-opinion" -> "Statement-opinion"
$ws.Range("I36").Value = "sv"
$ws.Range("J36").Value = "Statement-opinion"
